# feat: add 2022-Q1 data
#
#  - insert a new worksheet "2022-Q1" (per-fund holding detail for the
#    quarter) right before the existing "总计" (grand-total/summary) sheet
#  - prepend a 2022-Q1 row to the "总计" summary table, pushing the other
#    quarters down one row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet with per-fund holding detail, inserted just before 总计
# ---------------------------------------------------------------------------
$totalBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalBefore)
$q1.Name = "2022-Q1"

# NOTE: the worksheet reference used as the "insert before" anchor tracks the
# collection *position*, not the original sheet's identity, so after the
# insert it resolves to the brand-new sheet rather than "总计". Any further
# work on the summary sheet re-fetches it by name (see step 2).

# Header row (bold, centered, thin border - matches the other quarter sheets)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$header = $q1.Range("B1:H1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Force the code/ratio columns to stay text so leading/trailing zeros survive
# (e.g. fund code "010389", share "7.10") exactly like the other quarter sheets
$q1.Range("B2:G5").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "010389"
$q1.Range("C2").Value = "易方达科益混合A"
$q1.Range("D2").Value = "7.10"
$q1.Range("E2").Value = "92.94"
$q1.Range("F2").Value = "4.07"
$q1.Range("G2").Value = "0.2890"
$q1.Range("H2").Value = 8

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "011649"
$q1.Range("C3").Value = "易方达逆向投资混合A"
$q1.Range("D3").Value = "7.49"
$q1.Range("E3").Value = "85.02"
$q1.Range("F3").Value = "2.76"
$q1.Range("G3").Value = "0.2067"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "011650"
$q1.Range("C4").Value = "易方达逆向投资混合C"
$q1.Range("D4").Value = "1.96"
$q1.Range("E4").Value = "85.02"
$q1.Range("F4").Value = "2.76"
$q1.Range("G4").Value = "0.0541"
$q1.Range("H4").Value = 10

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "010390"
$q1.Range("C5").Value = "易方达科益混合C"
$q1.Range("D5").Value = "0.29"
$q1.Range("E5").Value = "92.94"
$q1.Range("F5").Value = "4.07"
$q1.Range("G5").Value = "0.0118"
$q1.Range("H5").Value = 8

# Column A (the row-index column) is bold/centered/top/bordered, like the
# other quarter sheets
$colA = $q1.Range("A2:A5")
$colA.Font.Bold = $true
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160
$colA.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2) Prepend a 2022-Q1 summary row into "总计"
# ---------------------------------------------------------------------------
# Re-fetch by name (see note above) - this is still the original summary sheet.
$total = $wb.Worksheets.Item("总计")

# Shift the six existing data rows down one slot by writing them into their
# new positions directly (row 7 is new, rows 2..6 already carry the correct
# per-column formatting from the original sheet, so overwriting their values
# in place - rather than doing a row Insert - keeps that formatting intact).
$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 2
$total.Range("D7").Value = 1.07

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 5
$total.Range("D6").Value = 4.68

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 4
$total.Range("D5").Value = 6.45

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.02

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 8
$total.Range("D3").Value = 3.29

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.5600000000000001

# Row 7 is brand new (the sheet previously ended at row 6), so it needs the
# same bold/centered/top/bordered look the other row-index cells (A2:A6) have
$a7 = $total.Range("A7")
$a7.Font.Bold = $true
$a7.HorizontalAlignment = -4108
$a7.VerticalAlignment = -4160
$a7.Borders.LineStyle = 1

# Restore the original active-sheet selection (unaffected by this edit)
$wb.Worksheets.Item(1).Activate()
